# Auto-generated edit script: apply updated currentAveragePrice / Leve profit figures
# per the "update Sheets via scheduled runner" commit.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H138").Value = 2313.974
$ws.Range("J138").Value = 2315.9138
$ws.Range("L138").Value = 6947.741399999999
$ws.Range("N138").Value = -17227.7414

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1012.2
$ws.Range("I2").Value = 1104.3334
$ws.Range("K2").Value = 1104.3334
$ws.Range("M2").Value = -991.3334
$ws.Range("H32").Value = 1055218.8
$ws.Range("I32").Value = 1478075.9
$ws.Range("J32").Value = 21568.223
$ws.Range("K32").Value = 1478075.9
$ws.Range("L32").Value = 21568.223
$ws.Range("M32").Value = -1477788.9
$ws.Range("N32").Value = -22142.223
$ws.Range("H74").Value = 1685.56
$ws.Range("I74").Value = 1551.8182
$ws.Range("K74").Value = 1551.8182
$ws.Range("M74").Value = -677.8181999999999
$ws.Range("H77").Value = 1685.56
$ws.Range("I77").Value = 1551.8182
$ws.Range("K77").Value = 7759.090999999999
$ws.Range("M77").Value = -3391.090999999999
$ws.Range("H110").Value = 2637
$ws.Range("I110").Value = 2637
$ws.Range("K110").Value = 2637
$ws.Range("M110").Value = -592
$ws.Range("H116").Value = 1012.2
$ws.Range("I116").Value = 1104.3334
$ws.Range("K116").Value = 1104.3334
$ws.Range("M116").Value = 1189.6666
$ws.Range("H132").Value = 2997.1667
$ws.Range("I132").Value = 2542.2942
$ws.Range("J132").Value = 4101.857
$ws.Range("K132").Value = 7626.882599999999
$ws.Range("L132").Value = 12305.571
$ws.Range("M132").Value = -5096.882599999999
$ws.Range("N132").Value = -17365.571

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1012.2
$ws.Range("I3").Value = 1104.3334
$ws.Range("K3").Value = 1104.3334
$ws.Range("M3").Value = -990.3334
$ws.Range("H134").Value = 2296.9656
$ws.Range("I134").Value = 2147.92
$ws.Range("J134").Value = 3228.5
$ws.Range("K134").Value = 6443.76
$ws.Range("L134").Value = 9685.5
$ws.Range("M134").Value = -3908.76
$ws.Range("N134").Value = -14755.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 1229.4286
$ws.Range("I58").Value = 910.6667
$ws.Range("J58").Value = 1468.5
$ws.Range("K58").Value = 910.6667
$ws.Range("L58").Value = 1468.5
$ws.Range("M58").Value = -707.6667
$ws.Range("N58").Value = -1874.5
$ws.Range("H68").Value = 23691.691
$ws.Range("J68").Value = 23691.691
$ws.Range("L68").Value = 23691.691
$ws.Range("N68").Value = -25189.691
$ws.Range("H71").Value = 23691.691
$ws.Range("J71").Value = 23691.691
$ws.Range("L71").Value = 71075.073
$ws.Range("N71").Value = -78563.073
$ws.Range("H134").Value = 1690.48
$ws.Range("I134").Value = 1407.7894
$ws.Range("J134").Value = 2585.6667
$ws.Range("K134").Value = 4223.3682
$ws.Range("L134").Value = 7757.000100000001
$ws.Range("M134").Value = -1688.3682
$ws.Range("N134").Value = -12827.0001
$ws.Range("H136").Value = 1229.4286
$ws.Range("I136").Value = 910.6667
$ws.Range("J136").Value = 1468.5
$ws.Range("K136").Value = 2732.0001
$ws.Range("L136").Value = 4405.5
$ws.Range("M136").Value = -182.0001000000002
$ws.Range("N136").Value = -9505.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H49").Value = 5415.6665
$ws.Range("I49").Value = 750
$ws.Range("J49").Value = 7748.5
$ws.Range("K49").Value = 2250
$ws.Range("L49").Value = 23245.5
$ws.Range("M49").Value = -2094
$ws.Range("N49").Value = -23557.5
$ws.Range("H68").Value = 1089.5264
$ws.Range("J68").Value = 1283.4166
$ws.Range("L68").Value = 3850.2498
$ws.Range("N68").Value = -5472.2498
$ws.Range("H71").Value = 1089.5264
$ws.Range("J71").Value = 1283.4166
$ws.Range("L71").Value = 11550.7494
$ws.Range("N71").Value = -19662.7494
$ws.Range("H92").Value = 798.3333
$ws.Range("I92").Value = 200
$ws.Range("J92").Value = 918
$ws.Range("K92").Value = 600
$ws.Range("L92").Value = 2754
$ws.Range("M92").Value = 648
$ws.Range("N92").Value = -5250
$ws.Range("H137").Value = 12292
$ws.Range("I137").Value = 15159.875
$ws.Range("J137").Value = 4644.3335
$ws.Range("K137").Value = 45479.625
$ws.Range("L137").Value = 13933.0005
$ws.Range("M137").Value = -40379.625
$ws.Range("N137").Value = -24133.0005
$ws.Range("H140").Value = 1511.92
$ws.Range("I140").Value = 1278.4706
$ws.Range("K140").Value = 3835.4118
$ws.Range("M140").Value = 1344.5882

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 2007
$ws.Range("I102").Value = 2000
$ws.Range("J102").Value = 2014
$ws.Range("K102").Value = 2000
$ws.Range("L102").Value = 2014
$ws.Range("M102").Value = -378
$ws.Range("N102").Value = -5258
$ws.Range("H132").Value = 2891.7036
$ws.Range("I132").Value = 2780.2354
$ws.Range("J132").Value = 3081.2
$ws.Range("K132").Value = 8340.706200000001
$ws.Range("L132").Value = 9243.599999999999
$ws.Range("M132").Value = -5810.706200000001
$ws.Range("N132").Value = -14303.6

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H63").Value = 30000
$ws.Range("J63").Value = 30000
$ws.Range("L63").Value = 30000
$ws.Range("N63").Value = -31498
$ws.Range("H66").Value = 30000
$ws.Range("J66").Value = 30000
$ws.Range("L66").Value = 90000
$ws.Range("N66").Value = -97488
$ws.Range("H68").Value = 2350.2183
$ws.Range("I68").Value = 1925.3334
$ws.Range("J68").Value = 2612.647
$ws.Range("K68").Value = 1925.3334
$ws.Range("L68").Value = 2612.647
$ws.Range("M68").Value = -1176.3334
$ws.Range("N68").Value = -4110.647
$ws.Range("H70").Value = 95000
$ws.Range("J70").Value = 95000
$ws.Range("L70").Value = 95000
$ws.Range("N70").Value = -95540
$ws.Range("H71").Value = 2350.2183
$ws.Range("I71").Value = 1925.3334
$ws.Range("J71").Value = 2612.647
$ws.Range("K71").Value = 9626.666999999999
$ws.Range("L71").Value = 13063.235
$ws.Range("M71").Value = -5882.666999999999
$ws.Range("N71").Value = -20551.235
$ws.Range("H73").Value = 95000
$ws.Range("J73").Value = 95000
$ws.Range("L73").Value = 95000
$ws.Range("N73").Value = -96872
$ws.Range("H132").Value = 3146.0417
$ws.Range("I132").Value = 1994.5
$ws.Range("K132").Value = 5983.5
$ws.Range("M132").Value = -3453.5
$ws.Range("H136").Value = 11113764
$ws.Range("I136").Value = 4460
$ws.Range("J136").Value = 16668417
$ws.Range("K136").Value = 13380
$ws.Range("L136").Value = 50005251
$ws.Range("M136").Value = -10830
$ws.Range("N136").Value = -50010351

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H126").Value = 1577.1111
$ws.Range("I126").Value = 1345.3334
$ws.Range("K126").Value = 4036.0002
$ws.Range("M126").Value = -1566.0002
$ws.Range("H136").Value = 2084.658
$ws.Range("I136").Value = 1923.7037
$ws.Range("J136").Value = 2479.7273
$ws.Range("K136").Value = 5771.1111
$ws.Range("L136").Value = 7439.1819
$ws.Range("M136").Value = -3221.1111
$ws.Range("N136").Value = -12539.1819

